$d = $word.ActiveDocument

# Title (appears twice, both change identically to the same new text)
$d.Content.Find.Execute(
    "Play Heroes’ Gathering Free - Exciting Features and Gameplay",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Heroes’ Gathering Free Slot - Exciting Gameplay & Creative Features",
    2)

# "What we like" bullet list
$d.Content.Find.Execute(
    "Unique board game-style design",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exciting gameplay",
    2)

$d.Content.Find.Execute(
    "Exciting gameplay with creative features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Creative features",
    2)

$d.Content.Find.Execute(
    "High max bet for potential big winnings",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unique board game design",
    2)

$d.Content.Find.Execute(
    "Above-average RTP",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Impressive graphics and animations",
    2)

# "What we don't like" bullet list
$d.Content.Find.Execute(
    "Mini-game winnings are smaller",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited number of paylines",
    2)

$d.Content.Find.Execute(
    "Only 20 paylines",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Mini-game winnings may be smaller",
    2)

# Meta description (italic) text at end of document
$d.Content.Find.Execute(
    "Experience the unique board game-style design and creative features of Heroes’ Gathering. Play free and enjoy exciting gameplay with high max bet and above-average RTP.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Heroes’ Gathering, a free online slot game with exciting gameplay and creative features.",
    2)
